$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The date strings in column A (rows 3-21) are stored as plain text.
# Some of them (day <= 12) are ambiguous and would otherwise be auto-
# parsed into real Excel dates when assigned, so force the range to
# Text format first, assign the new "DD-MM-YYYY" strings, then restore
# the original (default/"Normal") style so no formatting residue is left
# behind.
$dateRange = $ws.Range("A3:A21")
$dateRange.NumberFormat = "@"

$ws.Range("A3").Value = "28-07-2022"
$ws.Range("A4").Value = "01-08-2022"
$ws.Range("A5").Value = "04-08-2022"
$ws.Range("A6").Value = "08-08-2022"
$ws.Range("A7").Value = "11-08-2022"
$ws.Range("A8").Value = "15-08-2022"
$ws.Range("A9").Value = "18-08-2022"
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("A12").Value = "29-08-2022"
$ws.Range("A13").Value = "01-09-2022"
$ws.Range("A14").Value = "05-09-2022"
$ws.Range("A15").Value = "08-09-2022"
$ws.Range("A16").Value = "12-09-2022"
$ws.Range("A17").Value = "15-09-2022"
$ws.Range("A18").Value = "19-09-2022"
$ws.Range("A19").Value = "22-09-2022"
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("A21").Value = "29-09-2022"

$dateRange.Style = "Normal"

# Row 3 attendance counts: one real (D) and one invalid (G) record now
# recognized for this date.
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1
